$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 556, shifting existing rows 556-674 down to 557-675
$ws.Rows.Item(556).Insert()

# Populate the newly inserted row 556 with the new weekly price record
$ws.Range("A556").Value = 10
$ws.Range("B556").Value = "Vega Modelo de Temuco"
$ws.Range("C556").Value = "La Araucanía"
$ws.Range("D556").Value = 45258
$ws.Range("E556").Value = 9
$ws.Range("F556").Value = "Fruta"
$ws.Range("G556").Value = 100108
$ws.Range("H556").Value = "Tropicales y subtropicales"
$ws.Range("I556").Value = 100108002
$ws.Range("J556").Value = "Mango"
$ws.Range("K556").Value = "Sin especificar"
$ws.Range("L556").Value = "Primera"
$ws.Range("M556").Value = 215
$ws.Range("N556").Value = 13000
$ws.Range("O556").Value = 13000
$ws.Range("P556").Value = 13000
$ws.Range("Q556").Value = "$/bandeja 4 kilos"
$ws.Range("R556").Value = "Brasil"
$ws.Range("S556").Value = 3250
$ws.Range("T556").Value = 4
